$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project list")

# Update project data file name from project_1 to project_data_defaults
$ws.Range("B2").Value = "project_data_defaults"

# Number of turbines: 67 -> 100
$ws.Range("I2").Value = 100

# Distance to interconnect (miles): 5 -> 10
$ws.Range("W2").Value = 10

# Interconnect Voltage (kV): 130 -> 137
$ws.Range("X2").Value = 137

# Update the active selection on the sheet to C4
[void]$ws.Range("C4").Select()
